$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data: PERSONA, ATTIVITÀ, TEMPO (min), DATA
$ws.Range("A13").Value = "Luca"
$ws.Range("B13").Value = "GDPR"
$ws.Range("C13").Value = 195
$ws.Range("D13").Value = 43521

# Match style of the date column used by the other rows (D4:D12)
$ws.Range("D13").NumberFormat = "m/d/yy"

# Update the active selection to reflect the new cursor position (C14)
$ws.Range("C14").Select()
